# Update "想去人数" (want-to-go count) and related figures to match
# the regenerated data snapshot (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 114
$ws.Range("F4").Value = 111
$ws.Range("F5").Value = 293
$ws.Range("F6").Value = 180
$ws.Range("F7").Value = 1159
$ws.Range("F8").Value = 405
$ws.Range("F9").Value = 92
$ws.Range("F10").Value = 111
$ws.Range("F11").Value = 128
$ws.Range("F13").Value = 250
$ws.Range("F14").Value = 146
$ws.Range("F15").Value = 140
$ws.Range("F16").Value = 1331
$ws.Range("F17").Value = 502
$ws.Range("F18").Value = 190
$ws.Range("F19").Value = 299
$ws.Range("F21").Value = 676
$ws.Range("F22").Value = 1081
$ws.Range("F24").Value = 1934
$ws.Range("F25").Value = 2505
$ws.Range("F26").Value = 1272
$ws.Range("F27").Value = 55
$ws.Range("F28").Value = 198
$ws.Range("F30").Value = 838
$ws.Range("F31").Value = 759
$ws.Range("F32").Value = 1000
$ws.Range("F33").Value = 123
$ws.Range("F35").Value = 747
$ws.Range("F36").Value = 381
$ws.Range("F37").Value = 581
$ws.Range("F38").Value = 723
$ws.Range("F40").Value = 210

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 465
$ws.Range("G14").Value = 399

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 862

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 862
$ws.Range("F3").Value = 114
$ws.Range("F6").Value = 111
$ws.Range("F7").Value = 293
$ws.Range("F8").Value = 180
$ws.Range("F11").Value = 1159
$ws.Range("F12").Value = 405
$ws.Range("F13").Value = 92
$ws.Range("F14").Value = 111
$ws.Range("F16").Value = 128
$ws.Range("F17").Value = 250
$ws.Range("F19").Value = 146
$ws.Range("F20").Value = 140
$ws.Range("F21").Value = 1331
$ws.Range("F22").Value = 502
$ws.Range("F23").Value = 190
$ws.Range("F24").Value = 299
$ws.Range("F26").Value = 1082
$ws.Range("F27").Value = 2505
$ws.Range("F29").Value = 1272
$ws.Range("F30").Value = 55
$ws.Range("F34").Value = 198
$ws.Range("F36").Value = 838
$ws.Range("F39").Value = 759
$ws.Range("F40").Value = 1000
$ws.Range("F41").Value = 747
$ws.Range("F42").Value = 381
$ws.Range("F43").Value = 581
$ws.Range("F44").Value = 723
$ws.Range("F48").Value = 210

